# Update the cryptos listing with the latest scraped prices / 1h volume %
# change. Values in column D that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as text (matching the original
# European-style "thousands separated by dots" text strings, e.g.
# "63.032.84"), instead of silently re-interpreting them as numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.032.84'
$ws.Range("D3").Value = '2.472.33'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '''577.67'
$ws.Range("E5").Value = '  +0.62%  '
$ws.Range("D6").Value = '''146.90'
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("D9").Value = '2.470.66'
$ws.Range("E9").Value = '  +0.60%  '
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").Value = '''29.00'
$ws.Range("E14").Value = '  +6.55%  '
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '2.919.76'
$ws.Range("E16").Value = '  -1.45%  '
$ws.Range("D17").Value = '63.113.63'
$ws.Range("E17").Value = '  +0.49%  '
$ws.Range("D18").Value = '2.465.76'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = '''8.19'
$ws.Range("E19").Value = '  +3.58%  '
$ws.Range("D20").Value = '''11.03'
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("D21").Value = '''329.73'
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("E22").Value = '  +9.36%  '
$ws.Range("D24").Value = '''0.999'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").Value = '''66.27'
$ws.Range("D26").Value = '''665.85'
$ws.Range("E26").Value = '  +6.62%  '
$ws.Range("D27").Value = '''9.58'
$ws.Range("E27").Value = '  +14.17%  '
$ws.Range("E28").Value = '  +0.68%  '
$ws.Range("D29").Value = '2.591.86'
$ws.Range("E29").Value = '  +1.03%  '
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  +549.52%  '
$ws.Range("E31").Value = '  +2.44%  '
$ws.Range("D32").Value = '''8.08'
$ws.Range("E32").Value = '  -0.95%  '
$ws.Range("E33").Value = '  +0.97%  '
$ws.Range("E34").Value = '  -3.11%  '
$ws.Range("E35").Value = '  +3.56%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  +0.53%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").Value = '''5.44'
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '''152.42'
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("E40").Value = '  -0.44%  '
$ws.Range("E41").Value = '  +0.51%  '
$ws.Range("E42").Value = '  -0.44%  '
$ws.Range("D43").Value = '''1.76'
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("E45").Value = '  +8.05%  '
$ws.Range("D46").Value = '''151.69'
$ws.Range("E46").Value = '  +4.65%  '
$ws.Range("D47").Value = '''15.14'
$ws.Range("E47").Value = '  +25.05%  '
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("D49").Value = '''20.65'
$ws.Range("E49").Value = '  +2.18%  '
$ws.Range("E50").Value = '  +0.71%  '
$ws.Range("D51").Value = '''0.0513'
$ws.Range("E51").Value = '  -0.73%  '